# Automatische test-sync: 2025-06-22 21:53:50
# Adds a new "Offerte voor 500 stuks" entry to the Logs sheet and
# keeps the Dashboard category-count pivot in sync.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 49

$logs.Cells.Item($newRow, 1).Value = "Offerte voor 500 stuks"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$logs.Cells.Item($newRow, 4).Value = "Offerte / Prijsaanvraag"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,
Bedankt voor uw interesse in product X. Om een offerte voor 500 stuks te kunnen opstellen, hebben we meer informatie nodig. Kunt u de specificaties van het product en eventuele voorkeuren doorgeven, zoals kleur, maat, en eventuele extra functies? Nadat we deze gegevens hebben ontvangen, zullen we zo spoedig mogelijk een offerte voor u opstellen.
Met vriendelijke groet,
[Naam] - E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 21:53:21"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# --- 2. Update the "Dashboard" category pivot -------------------------------
# "Offerte / Prijsaanvraag" now outranks "Retour / Terugbetaling" (5 vs 5 ->
# offerte moves up to row 4, retour drops to row 5) and the offerte count
# increases from 4 to 5.
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(4, 2).Value = 5

$dash.Cells.Item(5, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(5, 2).Value = 5

# --- 3. Re-anchor the conditional formatting ranges to include row 49 ------
$dFormats = $logs.Range("D2:D48").FormatConditions
for ($i = 1; $i -le $dFormats.Count; $i++) {
    $dFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D49"))
}

$gFormats = $logs.Range("G2:G48").FormatConditions
for ($i = 1; $i -le $gFormats.Count; $i++) {
    $gFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G49"))
}
